# Apply updated cryptocurrency price/volume figures (cryptos list refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.893.95'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.44%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.563.32'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.44%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '598.97'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.81%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '135.29'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.40%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.562.12'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.42%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.493'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.38%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.122'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.41%  '

$ws.Range("E11").Value = '  -2.65%  '

$ws.Range("E12").Value = '  -0.25%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.167.54'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.44%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000181'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.10%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.568.10'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.50%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.92'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.24%  '

$ws.Range("E17").Value = '  +0.48%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '64.496.37'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.03%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.00'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.95%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.30'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.79%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.81'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.50%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '388.05'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.24%  '

$ws.Range("E23").Value = '  +3.90%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.708.29'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.48%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '73.89'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.79%  '

$ws.Range("E26").Value = '  +0.18%  '

$ws.Range("E27").Value = '  +3.43%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.69'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.33%  '

$ws.Range("E29").Value = '  +0.14%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.28'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.20%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.38'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.85%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.48'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +24.03%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.564.06'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.91%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.97'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.98%  '

$ws.Range("E36").Value = '  +0.35%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.90'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.50%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '168.52'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.57%  '

$ws.Range("E39").Value = '  +4.06%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.96'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.13%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0803'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.24%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.824'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.65%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '26.75'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.44%  '

$ws.Range("E44").Value = '  +0.22%  '

$ws.Range("E45").Value = '  -0.04%  '

$ws.Range("E46").Value = '  +2.07%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.20'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.64%  '

$ws.Range("E48").Value = '  +1.26%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.480.42'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +11.83%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.88'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.79%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.864'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +7.80%  '
